$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3261028230190277
$ws.Range("B1").Value = 2.454371690750122
$ws.Range("C1").Value = 8.756747245788574
$ws.Range("D1").Value = 2.097196102142334
$ws.Range("E1").Value = 1.178708076477051
